# Update F-column (想去人数) figures for rows 2-21 on both the "展览"
# and "全部类型" worksheets, and the G12 (最低票价) figure on both sheets.

$sheetNames = @("展览", "全部类型")

# Row -> new F value
$fUpdates = @{
    2  = 384
    3  = 998
    4  = 233
    5  = 1365
    6  = 8391
    7  = 50
    10 = 231
    11 = 141
    12 = 3366
    14 = 335
    15 = 48
    16 = 877
    18 = 1083
    20 = 141
    21 = 1934
}

# Row -> new G value
$gUpdates = @{
    12 = 50
}

foreach ($sheetName in $sheetNames) {
    $ws = $excel.ActiveWorkbook.Worksheets.Item($sheetName)

    foreach ($row in $fUpdates.Keys) {
        $ws.Range("F$row").Value = $fUpdates[$row]
    }

    foreach ($row in $gUpdates.Keys) {
        $ws.Range("G$row").Value = $gUpdates[$row]
    }
}
